$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values such as
# "5.20" or "0.600" keep their exact textual representation instead
# of being auto-converted to numbers (which would drop trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "62.067.31"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "2.443.07"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "583.61"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("D6").Value = "142.57"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").Value = "2.438.03"
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("E11").Value = "  +2.58%  "
$ws.Range("D12").Value = "5.20"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "0.341"
$ws.Range("E13").Value = "  -2.73%  "
$ws.Range("D14").Value = "26.31"
$ws.Range("E14").Value = "  -0.71%  "
$ws.Range("D15").Value = "0.0000175"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "2.865.93"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").Value = "61.985.10"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "2.430.43"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "10.72"
$ws.Range("E19").Value = "  -3.02%  "
$ws.Range("D20").Value = "7.19"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").Value = "325.92"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "4.09"
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "1.90"
$ws.Range("E24").Value = "  -5.71%  "
$ws.Range("D25").Value = "65.55"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").Value = "9.14"
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("D27").Value = "601.06"
$ws.Range("E27").Value = "  -4.46%  "
$ws.Range("D28").Value = "0.0₃0963"
$ws.Range("E28").Value = "  +0.50%  "
$ws.Range("E30").Value = "  +0.78%  "
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("D32").Value = "7.96"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").Value = "4.87"
$ws.Range("E35").Value = "  -2.20%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").Value = "1.42"
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("D38").Value = "0.375"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").Value = "152.84"
$ws.Range("E39").Value = "  +4.75%  "
$ws.Range("D40").Value = "18.37"
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("D41").Value = "5.26"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("E42").Value = "  +2.11%  "
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").Value = "141.66"
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("D47").Value = "3.61"
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("D48").Value = "0.0₆0264"
$ws.Range("E48").Value = "  +18.17%  "
$ws.Range("D49").Value = "0.600"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("D50").Value = "0.0519"
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("D51").Value = "19.78"
$ws.Range("E51").Value = "  +0.40%  "

# Clear the temporary text-number-format so the cells keep the same
# (default) style they originally had, matching the source workbook.
$ws.Range("D2:D51").ClearFormats()
